$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6
$ws.Range("G6").Value = 1.55
$ws.Range("H6").Value = 3.7
$ws.Range("S6").Value = 1.57
$ws.Range("T6").Value = 2.25
$ws.Range("Z6").Value = 10
$ws.Range("AE6").Value = 29
$ws.Range("AF6").Value = 126
$ws.Range("AI6").Value = 34
$ws.Range("AJ6").Value = 23
$ws.Range("AL6").Value = 67
$ws.Range("AM6").Value = 81
$ws.Range("AN6").Value = 3.2
$ws.Range("AT6").Value = 2.25
$ws.Range("BB6").Value = 301

# Row 7
$ws.Range("H7").Value = 2.88
$ws.Range("J7").Value = 3.2
$ws.Range("M7").Value = 1.13
$ws.Range("N7").Value = 6
$ws.Range("O7").Value = 1.62
$ws.Range("P7").Value = 2.2
$ws.Range("S7").Value = 1.67
$ws.Range("T7").Value = 2.1
$ws.Range("U7").Value = 2.38
$ws.Range("V7").Value = 1.53
$ws.Range("Y7").Value = 11
$ws.Range("AI7").Value = 17
$ws.Range("AP7").Value = 34
$ws.Range("AR7").Value = 101
$ws.Range("AT7").Value = 2.1
$ws.Range("AU7").Value = 10
$ws.Range("BB7").Value = 151

# Row 9
$ws.Range("G9").Value = 2.75
$ws.Range("I9").Value = 2.63
$ws.Range("J9").Value = 3.6
$ws.Range("L9").Value = 3.5
$ws.Range("AK9").Value = 26
$ws.Range("AN9").Value = 4.75
$ws.Range("AU9").Value = 9
$ws.Range("AX9").Value = 4.5
$ws.Range("Z9").Value = 29
$ws.Range("AA9").Value = 26

# Row 10
$ws.Range("N10").Value = 9

# Row 11
$ws.Range("G11").Value = 2.3
$ws.Range("H11").Value = 2.9
$ws.Range("I11").Value = 3.4
$ws.Range("J11").Value = 3.2
$ws.Range("S11").Value = 1.57
$ws.Range("T11").Value = 2.25
$ws.Range("X11").Value = 10
$ws.Range("AA11").Value = 23
$ws.Range("AG11").Value = 900
$ws.Range("AK11").Value = 34
$ws.Range("AN11").Value = 4.33
$ws.Range("AO11").Value = 15
$ws.Range("AT11").Value = 2.25

# Row 14
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = 3.4
$ws.Range("I14").Value = 3.4
$ws.Range("J14").Value = 2.75
$ws.Range("W14").Value = 7
$ws.Range("AD14").Value = 6.5
$ws.Range("AI14").Value = 17
$ws.Range("AL14").Value = 29

# Row 42
$ws.Range("G42").Value = 1.62
$ws.Range("H42").Value = 3.6
$ws.Range("I42").Value = 5.75
$ws.Range("J42").Value = 2.3
$ws.Range("M42").Value = 1.08
$ws.Range("N42").Value = 8
$ws.Range("O42").Value = 1.44
$ws.Range("P42").Value = 2.63
$ws.Range("Q42").Value = 2.4
$ws.Range("R42").Value = 1.53
$ws.Range("Z42").Value = 11
$ws.Range("AC42").Value = 7
$ws.Range("AD42").Value = 7.5
$ws.Range("AH42").Value = 11
$ws.Range("AI42").Value = 29
$ws.Range("AM42").Value = 67
$ws.Range("AO42").Value = 9
$ws.Range("AQ42").Value = 29
$ws.Range("AX42").Value = 7.5
